$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.1038603333333333
$ws.Range("H2").Value = 0.311581
$ws.Range("M2").Value = 1.594873333333333
$ws.Range("N2").Value = 4.78462
$ws.Range("O2").Value = 0.09372679355272211
$ws.Range("P2").Value = 0.09372679355272213
$ws.Range("Q2").Value = 0.1656440760244444
$ws.Range("R2").Value = 1.49079668422
$ws.Range("S2").Value = 0.09372679355272211
$ws.Range("T2").Value = 0.09372679355272213

$ws.Range("G3").Value = 0.1038603333333333
$ws.Range("H3").Value = 0.311581
$ws.Range("O3").Value = 0.2690834924840127
$ws.Range("P3").Value = 0.2690834924840128
$ws.Range("Q3").Value = 0.4755533054789999
$ws.Range("R3").Value = 4.279979749311
$ws.Range("S3").Value = 0.2690834924840127
$ws.Range("T3").Value = 0.2690834924840128

$ws.Range("G4").Value = 0.1038603333333333
$ws.Range("H4").Value = 0.311581
$ws.Range("M4").Value = 4.495828
$ws.Range("N4").Value = 13.487484
$ws.Range("O4").Value = 0.2642087832291055
$ws.Range("P4").Value = 0.2642087832291055
$ws.Range("Q4").Value = 0.4669381946893333
$ws.Range("R4").Value = 4.202443752204
$ws.Range("S4").Value = 0.2642087832291055
$ws.Range("T4").Value = 0.2642087832291055

$ws.Range("G5").Value = 0.1038603333333333
$ws.Range("H5").Value = 0.311581
$ws.Range("M5").Value = 6.346716
$ws.Range("N5").Value = 19.040148
$ws.Range("O5").Value = 0.3729809307341596
$ws.Range("P5").Value = 0.3729809307341597
$ws.Range("Q5").Value = 0.659172039332
$ws.Range("R5").Value = 5.932548353988
$ws.Range("S5").Value = 0.3729809307341596
$ws.Range("T5").Value = 0.3729809307341597
